# Updates crypto price/volume(1h) figures on the "cryptos" sheet to the
# latest scraped values (GitHub Actions symbol-list refresh).
# Each target D/E cell is plain text (e.g. "306.44", "1.09%"), so we
# assign with a leading apostrophe to force text entry (avoiding Excel's
# auto-conversion to a number/percentage), then reset the style to "Normal"
# so the cell keeps its original (default) formatting/style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'306.44"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'1.09%"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'35.97"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'0.89%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.019"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'-1.25%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.08078"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'0.23%"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.925"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-0.96%"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'2.31%"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'7.846"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'0.66%"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.9306"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'0.28%"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.1251"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'-17.31%"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.1916"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'1.10%"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.09210"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'2.22%"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'0.03517"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'1.95%"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.09927"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'0.66%"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'0.001432"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'1.10%"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'0.006651"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'14.40%"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'3.615"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'2.14%"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'3.068"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'2.59%"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'0.3438"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-0.12%"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'5.167"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'2.60%"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'-0.36%"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.2531"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'5.97%"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'0.04407"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'0.001235"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'2.33%"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.004730"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'-1.87%"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.0001300"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'6.08%"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'0.0003129"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'3.80%"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.01959"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'4.64%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.05180"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'7.93%"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.007545"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'2.98%"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.01010"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-4.51%"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.1373"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'1.94%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.002099"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'2.03%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.01070"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'10.15%"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00006392"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'2.75%"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.00000000750"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'0.34%"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'63.57"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-1.70%"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.001659"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'-0.01%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002100"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'0.34%"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.0002000"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'0.34%"
$c.Style = "Normal"
